$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - NroSiniestro value (text-formatted cell, keep quote-prefix style)
$ws.Range("F2").Formula = "'1120194100370"

# Row 3 - PREPROD environment now uses the "i-" (internal) host
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("F3").Formula = "'1120170200907"
$ws.Range("G3").Value = "Cheque"

# Row 4
$ws.Range("G4").Value = "Transferencia electrónica de fondos"

# Update sheet view: scroll so column B is the left-most visible column,
# and move the active selection to G7.
$ws.Activate()
$ws.Range("G7").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
